# Update pl_mw.xlsx result values for the 380 kV case (rows 2-25, columns C-N)
# Columns I and M remain 0 (unchanged) and are intentionally omitted.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ C=0.4067532613949538; D=0.07964120985750611; E=0.1662588220766459; F=3.194189273915754; G=2.425608670607772; H=1.969510745841774; J=0.2905240558401516; K=2.999162257781109; L=0.1420529065393481; N=1.791790579517773 }
    3 = @{ C=0.402842577632498; D=0.07826993896958356; E=0.1655160427446916; F=3.195487564848605; G=2.425051463014654; H=1.976853991231081; J=0.2904733534460249; K=2.85566803932619; L=0.1419392129316144; N=1.812790627531747 }
    4 = @{ C=0.400631243074983; D=0.07744739348875385; E=0.1651249368751628; F=3.198030570339768; G=2.426197032378354; H=1.982330368839271; J=0.29057192492251; K=2.768826768386361; L=0.1419129471612486; N=1.826364861297012 }
    5 = @{ C=0.399777855065139; D=0.07711712053670539; E=0.1649819135910811; F=3.199505216553149; G=2.427037072169469; H=1.984805106635193; J=0.290644701236765; K=2.733756823142357; L=0.1419132166564445; N=1.832067482746233 }
    6 = @{ C=0.3996390349336707; D=0.07706257729482502; E=0.1649591531147934; F=3.199776536752765; G=2.427199075538965; H=1.985230708845478; J=0.290658754920166; K=2.727952743531148; L=0.1419139248572066; N=1.833024727164901 }
    7 = @{ C=0.4006195406351338; D=0.07744291933687464; E=0.1651229417630233; F=3.198048683952919; G=2.426206851510727; H=1.982362760170986; J=0.2905727743909097; K=2.768352512201091; L=0.1419129063345892; N=1.826441076448315 }
    8 = @{ C=0.4053654604390999; D=0.07916438686101657; E=0.165989239770731; F=3.194274166834148; G=2.425107233183354; H=1.971841723584703; J=0.290479648493708; K=2.949423217748915; L=0.1420046816898548; N=1.79889023118896 }
    9 = @{ C=0.4161790320826526; D=0.0826928488063956; E=0.168202905810066; F=3.200758813786194; G=2.4347985329685; H=1.958899895660068; J=0.2913268620752802; K=3.314536202532679; L=0.1425292251775154; N=1.750259519736733 }
    10 = @{ C=0.4250446459775219; D=0.08537670665713648; E=0.1701427608778658; F=3.214041830435463; G=2.449206771620112; H=1.954099108236704; J=0.2925786317704322; K=3.588933063821116; L=0.1431236268240852; N=1.717819570718682 }
    11 = @{ C=0.4292784220224064; D=0.08661722428446694; E=0.1710932717717384; F=3.221946466217261; G=2.457358597949622; H=1.952941816086508; J=0.2932851238862213; K=3.715107656163809; L=0.1434392197761412; N=1.703775947069342 }
    12 = @{ C=0.4309105365254311; D=0.08708976512700417; E=0.1714629819780598; F=3.225208421741542; G=2.460676314157922; H=1.952651553183898; J=0.2935723863135422; K=3.763080942460476; L=0.1435652062370352; N=1.698560642482217 }
    13 = @{ C=0.4305577471100719; D=0.08698787178155243; E=0.1713829238703717; F=3.224493939332561; G=2.459951502439083; H=1.952707479398498; J=0.2935096414067573; K=3.752740428309608; L=0.1435377850063233; N=1.699679281622785 }
    14 = @{ C=0.4294121183186519; D=0.08665604492474444; E=0.1711234922826321; F=3.222209440249515; G=2.457626916687417; H=1.95291496822361; J=0.2933083616140451; K=3.719050567502222; L=0.1434494550425001; N=1.703344821655822 }
    15 = @{ C=0.4287141482397203; D=0.08645315282323196; E=0.1709658551667488; F=3.220845129727721; G=2.456233129759596; H=1.953061342901378; J=0.2931876418873074; K=3.698439768970843; L=0.1433961934341426; N=1.70560345030627 }
    16 = @{ C=0.424771999088307; D=0.08529602747234577; E=0.170082011044844; F=3.213562788662855; G=2.448706267800105; H=1.954195428042482; J=0.2925352205095777; K=3.580714386778368; L=0.1431039097542026; N=1.718751721015177 }
    17 = @{ C=0.4224050422092489; D=0.08459116638620401; E=0.1695572227024265; F=3.209572857648823; G=2.444498694302126; H=1.955154340940908; J=0.292170098804938; K=3.508839147731578; L=0.1429361641022311; N=1.72700057512602 }
    18 = @{ C=0.4210625275917153; D=0.08418759759548067; E=0.1692617854056735; F=3.207453173987318; G=2.442228902830806; H=1.955802490127439; J=0.2919729894925212; K=3.46762551146162; L=0.1428439365953977; N=1.731812245799926 }
    19 = @{ C=0.4206112205348802; D=0.08405127478719976; E=0.1691628562938412; F=3.206765553276455; G=2.441486170616457; H=1.956038524142144; J=0.2919084664734228; K=3.453693125697953; L=0.1428134414520059; N=1.733452924186736 }
    20 = @{ C=0.4226550530698887; D=0.08466600899286192; E=0.1696124243179611; F=3.209979451899258; G=2.444931035096403; H=1.95504226198463; J=0.2922076314872015; K=3.516477238507321; L=0.1429535806530353; N=1.726115520192216 }
    21 = @{ C=0.4297478335390963; D=0.08675343526075352; E=0.1711994285241971; F=3.222873155144157; G=2.458303431400594; H=1.952850004778071; J=0.2933669467527054; K=3.728940842430518; L=0.1434752240414028; N=1.702265375192386 }
    22 = @{ C=0.4345516876236957; D=0.08813389781563785; E=0.1722935818928555; F=3.232866198295994; G=2.468388837539493; H=1.952279939376751; J=0.294239633575053; K=3.86892698163615; L=0.1438538989312477; N=1.687276663882152 }
    23 = @{ C=0.4319723772935617; D=0.08739564803621391; E=0.1717044050322514; F=3.227389118539918; G=2.462882565968101; H=1.952505138061809; J=0.2937633338734855; K=3.794110585724184; L=0.1436483456501847; N=1.695221596551963 }
    24 = @{ C=0.4225419662212175; D=0.0846321674610877; E=0.1695874481350792; F=3.209795088143224; G=2.4447351091641; H=1.955092631199079; J=0.2921906230807494; K=3.513023719993328; L=0.1429456935101427; N=1.726515437684753 }
    25 = @{ C=0.4130921290563947; D=0.08172212646591248; E=0.1675489745471381; F=3.197512457633991; G=2.430901694581962; H=1.961575747047078; J=0.2909872657231745; K=3.214686139026355; L=0.1423505423340679; N=1.762837805269044 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
